$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "jr" label to rows 24-26 in column C (R type rows that are now jr rows)
$ws.Range("C24").Value = "jr"
$ws.Range("C25").Value = "jr"
$ws.Range("C26").Value = "jr"

# New column E (Control Signal names) -- entered in this order to match original authoring order
$ws.Range("E1").Value = "Control Signal"
$ws.Range("E2").Value = "branch"
$ws.Range("E4").Value = "PCwrt"
$ws.Range("E5").Value = "IRwrt"
$ws.Range("E6").Value = "memOWrt"
$ws.Range("E7").Value = "Awrt"
$ws.Range("E8").Value = "Bwrt"
$ws.Range("E9").Value = "ALUwrt"
$ws.Range("E10").Value = "regWrt"
$ws.Range("E11").Value = "memWrt"
$ws.Range("E12").Value = "wAdrs"
$ws.Range("E13").Value = "wDat[1:0]"
$ws.Range("E14").Value = "memAdrsSlct"
$ws.Range("E16").Value = "immSlct"
$ws.Range("E17").Value = "BNEoBEQ"
$ws.Range("E3").Value = "jump[1:0]"
$ws.Range("E15").Value = "imOrR[1:0]"

# New column F (descriptions) -- entered in this order to match original authoring order
$ws.Range("F1").Value = "What it does"
$ws.Range("F3").Value = "00 default, 01 if jumping, 10 if jumping to register"
$ws.Range("F4").Value = "high if writing to pc"
$ws.Range("F5").Value = "high if writing to IR"
$ws.Range("F6").Value = "high if writing to MemOWrt"
$ws.Range("F7").Value = "high if writing to AiA"
$ws.Range("F8").Value = "high if writing to AiB"
$ws.Range("F9").Value = "high if writing to ALUout"
$ws.Range("F10").Value = "high if writing to Reg. File"
$ws.Range("F11").Value = "high if writing to Mem File"
$ws.Range("F12").Value = "controls mux into wDest, default is `$m, high if address from IR"
$ws.Range("F13").Value = "controls mux into wData, 00 if from memO, 01 if from ALUout, 10 if from AiA, 11 if from Sign Ext."
$ws.Range("F14").Value = "controls mux into adrs of Mem, default pc, high if ALUout"
$ws.Range("F16").Value = "controls mux into Sign Ext., default is 12 bit immediate, high is 10 bit immediate"
$ws.Range("F17").Value = "controls mux out of isZero from ALU, default is isZero (used for BEQ), high is ~(isZero)  (used for BNE)"
$ws.Range("F15").Value = "controls mux into AiB, default r2out from Reg File, 01 from Sign Ext., 10 from PC+1 (for JR), 11 empty (currently)"
$ws.Range("F2").Value = "controls mux into pc adder, default is PC+1, high is Sign Ext. (for branching)"

# Column widths for E and F (best-fit to content, matching target bestFit columns)
$ws.Columns.Item(5).ColumnWidth = 12.67
$ws.Columns.Item(6).ColumnWidth = 100.67

# Update selection to C27 as in target
[void]$ws.Range("C27").Select()
